{"js": "// Apply the \"L\u1eddi c\u1ea3m \u01a1n\" wording fixes described by the commit diff:\n//  1. \"tr\u00e2n tr\u1ecdng, Ch\u00fang\" -> \"tr\u00e2n tr\u1ecdng, ch\u00fang\"\n//  2. \"d\u1ef1 \u00e1n; xin g\u1eedi\"    -> \"d\u1ef1 \u00e1n. Xin g\u1eedi\"\n//  3. \"C\u01a1 s\u1edf, Ph\u00f2ng T\u00e0i ch\u00ednh \u0111\u00e3\" -> \"C\u01a1 s\u1edf \u0111\u00e3\"\n//  4. \"v\u00e0 s\u1ef1 \u0111\u00e1nh gi\u00e1\"    -> \"v\u00e0 \u0111\u00e1nh gi\u00e1\"\n//  5. move the stray \"_GoBack\" bookmark from mid-paragraph 1 to right before\n//     the \"Nh\u00f3m sinh vi\u00ean th\u1ef1c hi\u1ec7n \u0111\u1ec1 t\u00e0i\" run at the end of the document.\n\nconst body = context.document.body;\n\nfunction replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  return context.sync().then(() => {\n    if (results.items.length === 0) {\n      throw new Error(`Text not found: ${searchText}`);\n    }\n    results.items[0].insertText(replacement, \"Replace\");\n    return context.sync();\n  });\n}\n\nawait replaceOnce(\"tr\u00e2n tr\u1ecdng, Ch\u00fang\", \"tr\u00e2n tr\u1ecdng, ch\u00fang\");\nawait replaceOnce(\"d\u1ef1 \u00e1n; xin g\u1eedi\", \"d\u1ef1 \u00e1n. Xin g\u1eedi\");\nawait replaceOnce(\"C\u01a1 s\u1edf, Ph\u00f2ng T\u00e0i ch\u00ednh \u0111\u00e3\", \"C\u01a1 s\u1edf \u0111\u00e3\");\nawait replaceOnce(\"v\u00e0 s\u1ef1 \u0111\u00e1nh gi\u00e1\", \"v\u00e0 \u0111\u00e1nh gi\u00e1\");\n\n// Relocate the \"_GoBack\" bookmark: delete wherever it currently sits, then\n// re-insert it immediately before the closing \"Nh\u00f3m sinh vi\u00ean th\u1ef1c hi\u1ec7n \u0111\u1ec1\n// t\u00e0i\" signature line.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst signature = body.search(\"Nh\u00f3m sinh vi\u00ean th\u1ef1c hi\u1ec7n \u0111\u1ec1 t\u00e0i\", { matchCase: true });\nsignature.load(\"items\");\nawait context.sync();\n\nif (signature.items.length > 0) {\n  const startRange = signature.items[0].getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Apply the \"L\u1eddi c\u1ea3m \u01a1n\" wording fixes described by the commit diff:\n#  1. \"tr\u00e2n tr\u1ecdng, Ch\u00fang\" -> \"tr\u00e2n tr\u1ecdng, ch\u00fang\"\n#  2. \"d\u1ef1 \u00e1n; xin g\u1eedi\"    -> \"d\u1ef1 \u00e1n. Xin g\u1eedi\"\n#  3. \"C\u01a1 s\u1edf, Ph\u00f2ng T\u00e0i ch\u00ednh \u0111\u00e3\" -> \"C\u01a1 s\u1edf \u0111\u00e3\"\n#  4. \"v\u00e0 s\u1ef1 \u0111\u00e1nh gi\u00e1\"    -> \"v\u00e0 \u0111\u00e1nh gi\u00e1\"\n#  5. move the stray \"_GoBack\" bookmark from mid-paragraph 1 to right before\n#     the \"Nh\u00f3m sinh vi\u00ean th\u1ef1c hi\u1ec7n \u0111\u1ec1 t\u00e0i\" run at the end of the document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-OnceInDoc($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\nReplace-OnceInDoc \"tr\u00e2n tr\u1ecdng, Ch\u00fang\" \"tr\u00e2n tr\u1ecdng, ch\u00fang\"\nReplace-OnceInDoc \"d\u1ef1 \u00e1n; xin g\u1eedi\" \"d\u1ef1 \u00e1n. Xin g\u1eedi\"\nReplace-OnceInDoc \"C\u01a1 s\u1edf, Ph\u00f2ng T\u00e0i ch\u00ednh \u0111\u00e3\" \"C\u01a1 s\u1edf \u0111\u00e3\"\nReplace-OnceInDoc \"v\u00e0 s\u1ef1 \u0111\u00e1nh gi\u00e1\" \"v\u00e0 \u0111\u00e1nh gi\u00e1\"\n\n# Relocate the \"_GoBack\" bookmark: delete wherever it currently sits, then\n# re-insert it immediately before the closing \"Nh\u00f3m sinh vi\u00ean th\u1ef1c hi\u1ec7n \u0111\u1ec1\n# t\u00e0i\" signature line.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Nh\u00f3m sinh vi\u00ean th\u1ef1c hi\u1ec7n \u0111\u1ec1 t\u00e0i\"\n$find2.Forward = $true\n$find2.Wrap = 0\n$found = $find2.Execute()\nif ($found) {\n    $sigRange = $find2.Parent\n    $startRange = $d.Range($sigRange.Start, $sigRange.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $startRange)\n}\n"}
